$d = $word.ActiveDocument

function ReplaceOnce($findText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $ok = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# 1. Split the "engagement via likes, shares, or comments (Karnowski et al., 2017; ...)" citation:
#    move the "cognitive involvement" framing to the Oledorf-Hirsch citation.
ReplaceOnce ", engagement via likes, shares, or comments (Karnowski et al., 2017; " ", engagement via cognitive involvement ("

# 2. After "(Oledorf-Hirsch, 2018)" add the "likes, shares, or comments (Karnowski et al., 2017)" clause,
#    then introduce the final clause with "and finally".
ReplaceOnce "), or political knowledge and participation (" ") as well as likes, shares, or comments (Karnowski et al., 2017), and finally political knowledge and participation ("

# 3. "dominant" -> "narrative"
ReplaceOnce "dominant frame" "narrative frame"

# 4. Insert "organizing " before "thinking in this area"
ReplaceOnce "thinking in this area" "organizing thinking in this area"

# 5. Remove the redundant "effects" word: "stratification effects via the 'Matthew" -> "stratification via the 'Matthew"
ReplaceOnce "), or stratification effects via the ‘Matthew " "), or stratification via the ‘Matthew "

# 6. Expand "depending on various socio-technical conditions." into the fuller sentence.
ReplaceOnce " equalization and stratification are similarly possible, depending on various socio-technical conditions." " equalization and stratification are similarly possible, depending on various levels of analysis, as well as the socio-technical conditions of platform use."

# 7. Add new content replacing the lone "Second, " paragraph text.
ReplaceOnce "Second, " "Despite these findings, digital inequalities persist because information flows are curated based on individual interests that favor information-rich social networks (Robinson et al., 2015), creating the possibility of so-called ‘social media news deserts’ (Barnidge & Xenos, 2021). That is, individuals exercise a considerable amount of agency over their information flows; the algorithms that filter content are anchored in user behaviors and preferences (Thorson et al., 2021). These systems take on an actuary dimension (DeVito, 2017) in that networks comprised of shared interests are pooled, enhancing pre-existing preferences for news and political information. "
